# "made R0sp = beta_ss/b"
#
# This edit:
#  - Renames the existing "Sheet1" to "PerDay" (values expressed on a
#    per-day basis, 365-day year).
#  - Duplicates that sheet into a new sheet "PerSeason" placed right after
#    it, where the same rate formulas are instead expressed over a 90-day
#    field season (1/90 in place of 1/365, matching the "per-season"
#    transform beta_ss/b).
#  - The G6 recovery-rate formula on PerDay is corrected to use 1/365
#    (it previously, incorrectly, used 1/33).
#  - The new PerSeason sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Duplicate the sheet so PerSeason starts as an exact copy of PerDay ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheets ---
$ws1.Name = "PerDay"
$ws2.Name = "PerSeason"

# --- PerDay: fix the recovery-rate formula to use a 365-day period ---
$ws1.Range("G6").Formula = "=1 - ((4+2+1+4)/(15+4+4+11))^((1/365))"

# --- PerSeason: re-express the daily-rate formulas over a 90-day season ---
$ws2.Range("E2").Formula = "=1-(0.11^(1/90))"
$ws2.Range("E5").Formula = "=1-(0.77^(1/90))"
$ws2.Range("E6").Formula = "=1-(((0.37+0.17)/2)^(1/90))"
$ws2.Range("G6").Formula = "=1 - ((4+2+1+4)/(15+4+4+11))^((1/90))"
$ws2.Range("E7").Formula = "=1-(0.0263^(1/90))"

# The duplicated sheet inherits the sortState of the original; PerSeason
# shouldn't carry that leftover sort metadata.
$ws2.Sort.SortFields.Clear() | Out-Null

# --- Selection / view state ---
$ws1.Range("G6").Select() | Out-Null
$ws2.Range("G6").Select() | Out-Null

# PerSeason (the 2nd tab) ends up the active sheet/tab.
$ws2.Activate() | Out-Null
